$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.922
$ws.Range("B3").Value = 6.798999999999999
$ws.Range("B5").Value = 6.834999999999999
$ws.Range("D7").Value = -7.172999999999999
$ws.Range("A9").Value = -21.56
$ws.Range("D9").Value = -7.688
$ws.Range("B11").Value = 6.93
$ws.Range("B12").Value = 6.582000000000001
$ws.Range("A13").Value = -21.832
$ws.Range("A16").Value = -20.918
$ws.Range("A18").Value = -21.868
$ws.Range("A20").Value = -21.605
$ws.Range("B21").Value = 6.923999999999999
$ws.Range("D21").Value = -7.600999999999999
